$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E width (matches target layout) ---
$ws.Columns.Item(5).ColumnWidth = 53

# --- Row 2 (Assistant Research Professor, first bullet) ---
$ws.Range("A2").Value = "Assistant Research Professor"
$ws.Range("B2").Style = "Normal"
$ws.Range("B2").Value = "2015-present"
$ws.Range("C2").Value = "Psychology Department, College of Education and Human Services"
$ws.Range("D2").Value = "Utah State University"
$ws.Range("E2").Value = "Director, The Statistical Consulting Studio"

# --- Row 3 (Assistant Research Professor, second bullet) ---
$ws.Range("A3").Value = "Assistant Research Professor"
$ws.Range("B3").Style = "Normal"
$ws.Range("B3").Value = "2015-present"
$ws.Range("C3").Value = "Psychology Department, College of Education and Human Services"
$ws.Range("D3").Value = "Utah State University"
$ws.Range("E3").Value = "Instructor, graduate student quantitatice methods and statistics courses"

# --- Row 4 (Statistician, Office of Research Services, first bullet) ---
$ws.Range("A4").Value = "Statistician"
$ws.Range("B4").Value = "2013-2015"
$ws.Range("C4").Value = "Office of Research Services, College of Education and Human Services"
$ws.Range("D4").Value = "Utah State University"
$ws.Range("E4").Value = "Acting Director, Office of Methodological and Data Sciences"

# --- Row 5 (Statistician, Office of Research Services, second bullet) ---
$ws.Range("A5").Value = "Statistician"
$ws.Range("B5").Value = "2013-2015"
$ws.Range("C5").Value = "Office of Research Services, College of Education and Human Services"
$ws.Range("D5").Value = "Utah State University"
$ws.Range("E5").Value = "Instructor, graduate student quantitatice methods and statistics courses"

# --- Row 6 (Statistician and Data Manager, Center for Epidemiology, 1st bullet) ---
$ws.Range("A6").Value = "Statistician and Data Manager"
$ws.Range("B6").Value = "2005-2013"
$ws.Range("C6").Value = "Center for Epidemiology"
$ws.Range("D6").Value = "Utah State University"
$ws.Range("E6").Value = "Managed databases, clean data, and prepare custom datasets"

# --- Row 7 (Statistician and Data Manager, 2nd bullet, wrapped text) ---
$ws.Range("A7").Value = "Statistician and Data Manager"
$ws.Range("B7").Value = "2005-2013"
$ws.Range("C7").Value = "Center for Epidemiology"
$ws.Range("D7").Value = "Utah State University"
$ws.Range("E7").Value = "Performed statistical analyses and prepared publications, posters, presentations,`nand grant submissions"
$ws.Range("E7").WrapText = $true

# --- Row 8 (Statistician and Data Manager, 3rd bullet, wrapped text) ---
$ws.Range("A8").Value = "Statistician and Data Manager"
$ws.Range("B8").Value = "2005-2013"
$ws.Range("C8").Value = "Center for Epidemiology"
$ws.Range("D8").Value = "Utah State University"
$ws.Range("E8").Value = "Worked under three main grants and many co-investigators: University of Utah,`n  BYU, Duke, John Hopkins, University of Maryland, ect."
$ws.Range("E8").WrapText = $true

# --- Row 9 (Data Manager, Kenoi Genetics Lab) ---
$ws.Range("A9").Value = "Data Manager"
$ws.Range("B9").Value = 2012
$ws.Range("C9").Value = "Kenoi Genetics Lab"
$ws.Range("D9").Value = "Brigham Young University"
$ws.Range("E9").Value = "Managed databases and prepare custom datasets"

# --- Row 10 (Adjunct Lecturer) ---
$ws.Range("A10").Value = "Adjunct Lecturer"
$ws.Range("B10").Value = "2006-2008"
$ws.Range("C10").Value = "Mathematics and Statistics Department"
$ws.Range("D10").Value = "Utah State University"
$ws.Range("E10").Value = "Traditional, evening, and distance courses"

# --- Row 11 (High School Teacher, Math and Science, Sky View High School) ---
$ws.Range("A11").Value = "High School Teacher, Math and Science"
$ws.Range("B11").Value = "2000-2004"
$ws.Range("E11").Value = "Sky View High School, Smithfield, Utah"

# --- Row 12 (High School Teacher, Math and Science, Logan River Academy) ---
$ws.Range("A12").Value = "High School Teacher, Math and Science"
$ws.Range("B12").Value = "2000-2004"
$ws.Range("E12").Value = "Logan River Academy, Logan, Utah"

# --- Sheet-level formatting to match target layout ---
$ws.Rows.Item(1).RowHeight = 16.5
$ws.Range("A1:E12").RowHeight = 16.5

[void]$ws.Range("C13").Select()
